# Neutralize racial-misclassification language in this non-electoral
# (data_engineering) resume: "affecting all Black and Asian-American
# voters" -> "affecting 50M voters" (plus "nationwide" in the projects
# Impact line), per the commit message. Three spots are touched:
#   1. PROFESSIONAL SUMMARY paragraph (plain text swap)
#   2. Siege Analytics bullet point (the new "50M" becomes its own
#      bold, dark-slate-colored run, matching the styling already used
#      for the "23%"/"64%" figures in that same bullet)
#   3. KEY PROJECTS "Impact:" line for the classification system (plain
#      text swap, also adds "nationwide")

$d = $word.ActiveDocument

# Find the single paragraph whose text contains $marker.
function Get-ParagraphContaining($doc, $marker) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "*$marker*") {
            return $p
        }
    }
    return $null
}

# --- 1. Professional summary paragraph ---------------------------------
$pSummary = Get-ParagraphContaining $d "Data engineering professional with 15+ years"
if ($pSummary -ne $null) {
    $null = $pSummary.Range.Find.Execute(
        "all Black and Asian-American voters", $true, $false, $false,
        $false, $false, $true, 1, $false, "50M voters", 2)
}

# --- 2. Siege Analytics bullet: "50M" gets its own bold/colored run ----
$pBullet = Get-ParagraphContaining $d "Discovered systematic race coding errors"
if ($pBullet -ne $null) {
    $rngBullet = $pBullet.Range
    $null = $rngBullet.Find.Execute(
        "all Black and Asian-American", $true, $false, $false,
        $false, $false, $true, 1, $false, "50M", 2)
    # $rngBullet now collapses onto the just-inserted "50M" text;
    # style it like the other bolded metrics in this bullet (bold,
    # font color 2C3E50 -> RGB(44,62,80) -> Word BGR int 5258796).
    $rngBullet.Font.Bold = 1
    $rngBullet.Font.Color = 5258796
}

# --- 3. Key Projects "Impact:" line -------------------------------------
$pImpact = Get-ParagraphContaining $d "Impact: Corrected demographic data affecting"
if ($pImpact -ne $null) {
    $null = $pImpact.Range.Find.Execute(
        "all Black and Asian-American voters", $true, $false, $false,
        $false, $false, $true, 1, $false, "50M voters nationwide", 2)
}
